$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 175, shifting existing rows
# 175-194 down to 177-196 (preserving all of their data/formatting).
$ws.Rows(175).Resize(2).Insert()

# Populate the two newly inserted rows (175-176) with their data.
$ws.Range("A175").Value2 = 11
$ws.Range("B175").Value2 = "Vega Monumental Concepción"
$ws.Range("C175").Value2 = "Bíobío"
$ws.Range("D175").Value2 = 44783
$ws.Range("E175").Value2 = 8
$ws.Range("F175").Value2 = "Fruta"
$ws.Range("G175").Value2 = 100101
$ws.Range("H175").Value2 = "Berries"
$ws.Range("I175").Value2 = 100101007
$ws.Range("J175").Value2 = "Kiwi"
$ws.Range("K175").Value2 = "Hayward"
$ws.Range("L175").Value2 = "Especial"
$ws.Range("M175").Value2 = 250
$ws.Range("N175").Value2 = 8000
$ws.Range("O175").Value2 = 8000
$ws.Range("P175").Value2 = 8000
$ws.Range("Q175").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R175").Value2 = "Región de O'Higgins"
$ws.Range("S175").Value2 = 444
$ws.Range("T175").Value2 = 18

$ws.Range("A176").Value2 = 11
$ws.Range("B176").Value2 = "Vega Monumental Concepción"
$ws.Range("C176").Value2 = "Bíobío"
$ws.Range("D176").Value2 = 44783
$ws.Range("E176").Value2 = 8
$ws.Range("F176").Value2 = "Fruta"
$ws.Range("G176").Value2 = 100101
$ws.Range("H176").Value2 = "Berries"
$ws.Range("I176").Value2 = 100101007
$ws.Range("J176").Value2 = "Kiwi"
$ws.Range("K176").Value2 = "Hayward"
$ws.Range("L176").Value2 = "Primera"
$ws.Range("M176").Value2 = 550
$ws.Range("N176").Value2 = 6500
$ws.Range("O176").Value2 = 7000
$ws.Range("P176").Value2 = 6818
$ws.Range("Q176").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R176").Value2 = "Región de O'Higgins"
$ws.Range("S176").Value2 = 379
$ws.Range("T176").Value2 = 18
